# Update progress ("voortgang") for Koen & Djordy's tasks in the sprint 1 planning sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E6").Value = 1
